$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '84.727.99'
$ws.Range('E2').Value = '  +6.37%  '
$ws.Range('D3').Value = '3.245.85'
$ws.Range('E3').Value = '  +1.34%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = '''218.09'
$ws.Range('E5').Value = '  +3.22%  '
$ws.Range('D6').Value = '''627.36'
$ws.Range('E6').Value = '  -1.85%  '
$ws.Range('D7').Value = '''0.322'
$ws.Range('E7').Value = '  +30.80%  '
$ws.Range('E8').Value = '  -0.08%  '
$ws.Range('D9').Value = '''0.592'
$ws.Range('E9').Value = '  -1.54%  '
$ws.Range('D10').Value = '3.245.03'
$ws.Range('E10').Value = '  +1.31%  '
$ws.Range('E11').Value = '  -0.61%  '
$ws.Range('D12').Value = '''0.0000282'
$ws.Range('E12').Value = '  +10.21%  '
$ws.Range('D13').Value = '''0.166'
$ws.Range('E13').Value = '  -0.17%  '
$ws.Range('B14').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C14').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D14').Value = '3.843.36'
$ws.Range('E14').Value = '  +1.39%  '
$ws.Range('B15').Value = 'Toncoin'
$ws.Range('C15').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D15').Value = '''5.38'
$ws.Range('E15').Value = '  -0.56%  '
$ws.Range('D16').Value = '''32.79'
$ws.Range('E16').Value = '  +2.19%  '
$ws.Range('D17').Value = '84.373.98'
$ws.Range('E17').Value = '  +6.17%  '
$ws.Range('D18').Value = '3.240.05'
$ws.Range('E18').Value = '  +1.48%  '
$ws.Range('E19').Value = '  +7.92%  '
$ws.Range('D20').Value = '''14.37'
$ws.Range('E20').Value = '  -1.60%  '
$ws.Range('D21').Value = '''449.44'
$ws.Range('E21').Value = '  +2.31%  '
$ws.Range('D22').Value = '''9.13'
$ws.Range('E22').Value = '  -2.48%  '
$ws.Range('D23').Value = '''5.20'
$ws.Range('E23').Value = '  -1.30%  '
$ws.Range('E24').Value = '  +8.13%  '
$ws.Range('D25').Value = '''5.19'
$ws.Range('E25').Value = '  +7.59%  '
$ws.Range('D26').Value = '''11.93'
$ws.Range('E26').Value = '  +9.54%  '
$ws.Range('D27').Value = '3.394.42'
$ws.Range('E27').Value = '  +0.75%  '
$ws.Range('D28').Value = '''78.53'
$ws.Range('E28').Value = '  +1.42%  '
$ws.Range('E29').Value = '  +0.18%  '
$ws.Range('D30').Value = '''0.0000125'
$ws.Range('E30').Value = '  +1.83%  '
$ws.Range('D31').Value = '''9.16'
$ws.Range('E31').Value = '  +0.20%  '
$ws.Range('D32').Value = '''0.997'
$ws.Range('E32').Value = '  +0.12%  '
$ws.Range('D33').Value = '''0.156'
$ws.Range('E33').Value = '  +27.66%  '
$ws.Range('D34').Value = '''567.02'
$ws.Range('E34').Value = '  +0.84%  '
$ws.Range('E35').Value = '  -2.27%  '
$ws.Range('D36').Value = '''0.154'
$ws.Range('E36').Value = '  -2.79%  '
$ws.Range('D37').Value = '''2.00'
$ws.Range('E37').Value = '  -2.42%  '
$ws.Range('D38').Value = '''23.20'
$ws.Range('E38').Value = '  +0.63%  '
$ws.Range('B39').Value = 'FirstDigitalUSD'
$ws.Range('C39').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D39').Value = '''0.999'
$ws.Range('E39').Value = '  -0.05%  '
$ws.Range('B40').Value = 'RenderToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D40').Value = '''6.16'
$ws.Range('E40').Value = '  +8.60%  '
$ws.Range('D41').Value = '''0.407'
$ws.Range('E41').Value = '  -0.95%  '
$ws.Range('B42').Value = 'Stacks'
$ws.Range('C42').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D42').Value = '''2.05'
$ws.Range('E42').Value = '  +12.50%  '
$ws.Range('B43').Value = 'dogwifhat'
$ws.Range('C43').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D43').Value = '''3.08'
$ws.Range('E43').Value = '  +15.29%  '
$ws.Range('B44').Value = 'WhiteBITCoin'
$ws.Range('C44').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D44').Value = '''20.95'
$ws.Range('E44').Value = '  +3.42%  '
$ws.Range('D45').Value = '''160.20'
$ws.Range('E45').Value = '  -1.94%  '
$ws.Range('E46').Value = '  +0.07%  '
$ws.Range('D47').Value = '''188.39'
$ws.Range('E47').Value = '  -1.93%  '
$ws.Range('D48').Value = '''44.75'
$ws.Range('E48').Value = '  +3.94%  '
$ws.Range('E49').Value = '  -1.09%  '
$ws.Range('D50').Value = '''0.781'
$ws.Range('E50').Value = '  -2.55%  '
$ws.Range('D51').Value = '''25.97'
$ws.Range('E51').Value = '  +0.59%  '
